$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M94").ClearContents()
$ws.Range("H2").Value = 14577915
$ws.Range("J2").Value = 38461584
$ws.Range("L2").Value = 38461584
$ws.Range("N2").Value = -38461810
$ws.Range("H12").Value = 83.333336
$ws.Range("I12").Value = 55
$ws.Range("J12").Value = 140
$ws.Range("K12").Value = 55
$ws.Range("L12").Value = 140
$ws.Range("M12").Value = 115
$ws.Range("N12").Value = -480
$ws.Range("H18").Value = 620.55884
$ws.Range("I18").Value = 325.7742
$ws.Range("J18").Value = 3666.6667
$ws.Range("K18").Value = 325.7742
$ws.Range("L18").Value = 3666.6667
$ws.Range("M18").Value = -41.77420000000001
$ws.Range("N18").Value = -4234.6667
$ws.Range("H33").Value = 202.8108
$ws.Range("I33").Value = 194.38095
$ws.Range("J33").Value = 213.875
$ws.Range("K33").Value = 194.38095
$ws.Range("L33").Value = 213.875
$ws.Range("M33").Value = 34.61904999999999
$ws.Range("N33").Value = -671.875
$ws.Range("I40").Value = 6251580
$ws.Range("J40").Value = 3006
$ws.Range("K40").Value = 6251580
$ws.Range("L40").Value = 3006
$ws.Range("M40").Value = -6251405
$ws.Range("N40").Value = -3356
$ws.Range("H43").Value = 19232776
$ws.Range("I43").Value = 45456456
$ws.Range("J43").Value = 2079.6667
$ws.Range("K43").Value = 45456456
$ws.Range("L43").Value = 2079.6667
$ws.Range("M43").Value = -45456387
$ws.Range("N43").Value = -2217.6667
$ws.Range("H62").Value = 14715462
$ws.Range("I62").Value = 25013700
$ws.Range("J62").Value = 3692.8572
$ws.Range("K62").Value = 25013700
$ws.Range("L62").Value = 3692.8572
$ws.Range("M62").Value = -25013076
$ws.Range("N62").Value = -4940.8572
$ws.Range("H64").Value = 3090.9707
$ws.Range("I64").Value = 2822.4092
$ws.Range("J64").Value = 3583.3333
$ws.Range("K64").Value = 2822.4092
$ws.Range("L64").Value = 3583.3333
$ws.Range("M64").Value = -2574.4092
$ws.Range("N64").Value = -4079.3333
$ws.Range("H65").Value = 14715462
$ws.Range("I65").Value = 25013700
$ws.Range("J65").Value = 3692.8572
$ws.Range("K65").Value = 125068500
$ws.Range("L65").Value = 18464.286
$ws.Range("M65").Value = -125065380
$ws.Range("N65").Value = -24704.286
$ws.Range("H67").Value = 3090.9707
$ws.Range("I67").Value = 2822.4092
$ws.Range("J67").Value = 3583.3333
$ws.Range("K67").Value = 2822.4092
$ws.Range("L67").Value = 3583.3333
$ws.Range("M67").Value = -1964.4092
$ws.Range("N67").Value = -5299.3333
$ws.Range("H70").Value = 1973.4445
$ws.Range("I70").Value = 2792.2
$ws.Range("J70").Value = 950
$ws.Range("K70").Value = 8376.599999999999
$ws.Range("L70").Value = 2850
$ws.Range("M70").Value = -8106.599999999999
$ws.Range("N70").Value = -3390
$ws.Range("H73").Value = 1973.4445
$ws.Range("I73").Value = 2792.2
$ws.Range("J73").Value = 950
$ws.Range("K73").Value = 8376.599999999999
$ws.Range("L73").Value = 2850
$ws.Range("M73").Value = -7440.599999999999
$ws.Range("N73").Value = -4722
$ws.Range("H88").Value = 18255598
$ws.Range("J88").Value = 20283664
$ws.Range("L88").Value = 20283664
$ws.Range("N88").Value = -20284476
$ws.Range("H91").Value = 18255598
$ws.Range("J91").Value = 20283664
$ws.Range("L91").Value = 20283664
$ws.Range("N91").Value = -20286472
$ws.Range("H94").Value = 2368.3333
$ws.Range("I94").Value = 2368.3333
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2368.3333
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = -1917.3333
$ws.Range("H97").Value = 67000936
$ws.Range("I97").Value = 450
$ws.Range("J97").Value = 71786690
$ws.Range("K97").Value = 1350
$ws.Range("L97").Value = 215360070
$ws.Range("M97").Value = -854
$ws.Range("N97").Value = -215361062
$ws.Range("H99").Value = 292.62857
$ws.Range("I99").Value = 260.36365
$ws.Range("J99").Value = 825
$ws.Range("K99").Value = 781.09095
$ws.Range("L99").Value = 2475
$ws.Range("M99").Value = 716.90905
$ws.Range("N99").Value = -5471
$ws.Range("H100").Value = 8406
$ws.Range("I100").Value = 12121.111
$ws.Range("J100").Value = 2833.3333
$ws.Range("K100").Value = 12121.111
$ws.Range("L100").Value = 2833.3333
$ws.Range("M100").Value = -11580.111
$ws.Range("N100").Value = -3915.3333
$ws.Range("H103").Value = 557.1539
$ws.Range("J103").Value = 510
$ws.Range("L103").Value = 1530
$ws.Range("N103").Value = -2702
$ws.Range("H106").Value = 111112936
$ws.Range("J106").Value = 3366.6667
$ws.Range("L106").Value = 3366.6667
$ws.Range("N106").Value = -4628.6667
$ws.Range("H113").Value = 3550.842
$ws.Range("J113").Value = 3466.4
$ws.Range("L113").Value = 3466.4
$ws.Range("N113").Value = -9974.4
$ws.Range("H129").Value = 908.9595
$ws.Range("I129").Value = 914.53845
$ws.Range("J129").Value = 907.7705
$ws.Range("K129").Value = 2743.61535
$ws.Range("L129").Value = 2723.3115
$ws.Range("M129").Value = 2256.38465
$ws.Range("N129").Value = -12723.3115
$ws.Range("H137").Value = 1361.4286
$ws.Range("I137").Value = 935.11536
$ws.Range("J137").Value = 2593
$ws.Range("K137").Value = 2805.34608
$ws.Range("L137").Value = 7779
$ws.Range("M137").Value = -255.3460800000003
$ws.Range("N137").Value = -12879
$ws.Range("H138").Value = 2696.3076
$ws.Range("I138").Value = 1367.92
$ws.Range("J138").Value = 3322.9058
$ws.Range("K138").Value = 4103.76
$ws.Range("L138").Value = 9968.7174
$ws.Range("M138").Value = 1036.24
$ws.Range("N138").Value = -20248.7174

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2139.8
$ws.Range("I63").Value = 2066.5557
$ws.Range("K63").Value = 2066.5557
$ws.Range("M63").Value = -1380.5557
$ws.Range("H66").Value = 2139.8
$ws.Range("I66").Value = 2066.5557
$ws.Range("K66").Value = 10332.7785
$ws.Range("M66").Value = -6900.7785

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 15152181
$ws.Range("I107").Value = 20833860
$ws.Range("J107").Value = 1036.5555
$ws.Range("K107").Value = 20833860
$ws.Range("L107").Value = 1036.5555
$ws.Range("M107").Value = -20831940
$ws.Range("N107").Value = -4876.5555

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3484.4524
$ws.Range("I31").Value = 2435
$ws.Range("J31").Value = 7074.684
$ws.Range("K31").Value = 2435
$ws.Range("L31").Value = 7074.684
$ws.Range("M31").Value = -2140
$ws.Range("N31").Value = -7664.684
$ws.Range("H34").Value = 3484.4524
$ws.Range("I34").Value = 2435
$ws.Range("J34").Value = 7074.684
$ws.Range("K34").Value = 2435
$ws.Range("L34").Value = 7074.684
$ws.Range("M34").Value = -2233
$ws.Range("N34").Value = -7478.684
$ws.Range("H58").Value = 35715360
$ws.Range("I58").Value = 90910250
$ws.Range("J58").Value = 1022
$ws.Range("K58").Value = 90910250
$ws.Range("L58").Value = 1022
$ws.Range("M58").Value = -90910047
$ws.Range("N58").Value = -1428
$ws.Range("H62").Value = 4387.9
$ws.Range("J62").Value = 4747.5
$ws.Range("L62").Value = 4747.5
$ws.Range("N62").Value = -5995.5
$ws.Range("H65").Value = 4387.9
$ws.Range("J65").Value = 4747.5
$ws.Range("L65").Value = 23737.5
$ws.Range("N65").Value = -29977.5
$ws.Range("H136").Value = 35715360
$ws.Range("I136").Value = 90910250
$ws.Range("J136").Value = 1022
$ws.Range("K136").Value = 272730750
$ws.Range("L136").Value = 3066
$ws.Range("M136").Value = -272728200
$ws.Range("N136").Value = -8166

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 1200
$ws.Range("I94").Value = 300
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 9000
$ws.Range("M94").Value = -224
$ws.Range("N94").Value = -10352
$ws.Range("H131").Value = 696.9299999999999
$ws.Range("J131").Value = 775.2875
$ws.Range("L131").Value = 2325.8625
$ws.Range("N131").Value = -12405.8625
$ws.Range("H139").Value = 334099.38
$ws.Range("I139").Value = 1973.75
$ws.Range("J139").Value = 666225
$ws.Range("K139").Value = 5921.25
$ws.Range("L139").Value = 1998675
$ws.Range("M139").Value = -781.25
$ws.Range("N139").Value = -2008955

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 39147.168
$ws.Range("J42").Value = 39147.168
$ws.Range("L42").Value = 39147.168
$ws.Range("N42").Value = -40117.168
$ws.Range("H115").Value = 39147.168
$ws.Range("J115").Value = 39147.168
$ws.Range("L115").Value = 39147.168
$ws.Range("N115").Value = -41497.168
$ws.Range("H132").Value = 10113.5
$ws.Range("I132").Value = 712.5
$ws.Range("J132").Value = 14814
$ws.Range("K132").Value = 2137.5
$ws.Range("L132").Value = 44442
$ws.Range("M132").Value = 392.5
$ws.Range("N132").Value = -49502

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4168834.8
$ws.Range("I46").Value = 6945061
$ws.Range("J46").Value = 4495.5
$ws.Range("K46").Value = 6945061
$ws.Range("L46").Value = 4495.5
$ws.Range("M46").Value = -6944873
$ws.Range("N46").Value = -4871.5
$ws.Range("H68").Value = 1504.8422
$ws.Range("I68").Value = 1506.5714
$ws.Range("K68").Value = 1506.5714
$ws.Range("M68").Value = -757.5714
$ws.Range("H71").Value = 1504.8422
$ws.Range("I71").Value = 1506.5714
$ws.Range("K71").Value = 7532.857
$ws.Range("M71").Value = -3788.857
$ws.Range("H93").Value = 1299.6
$ws.Range("I93").Value = 1235.8182
$ws.Range("J93").Value = 1475
$ws.Range("K93").Value = 1235.8182
$ws.Range("L93").Value = 1475
$ws.Range("M93").Value = 12.18180000000007
$ws.Range("N93").Value = -3971
$ws.Range("H132").Value = 18872806
$ws.Range("I132").Value = 38463610
$ws.Range("J132").Value = 7589
$ws.Range("K132").Value = 115390830
$ws.Range("L132").Value = 22767
$ws.Range("M132").Value = -115388300
$ws.Range("N132").Value = -27827

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 30130
$ws.Range("J75").Value = 30130
$ws.Range("L75").Value = 30130
$ws.Range("N75").Value = -32002
$ws.Range("H78").Value = 30130
$ws.Range("J78").Value = 30130
$ws.Range("L78").Value = 90390
$ws.Range("N78").Value = -99750
$ws.Range("I113").Value = 71429064
$ws.Range("K113").Value = 214287192
$ws.Range("M113").Value = -214285022
